$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Tip" column (C) to match the Home_Team (A) pick for these rows
$ws.Range("C3").Value = "Roosters"
$ws.Range("C4").Value = "Rabbitohs"
$ws.Range("C6").Value = "Wests Tigers"
$ws.Range("C8").Value = "Panthers"
$ws.Range("C9").Value = "Storm"
